# Rename sheets cor_k7 -> cor_k07 and cor_k9 -> cor_k09, and update the
# corresponding labels/descriptions on the "Key" sheet to match.

$wb = $excel.ActiveWorkbook

# Rename the worksheet tabs.
$wb.Worksheets.Item("cor_k7").Name = "cor_k07"
$wb.Worksheets.Item("cor_k9").Name = "cor_k09"

# Update the "Key" sheet rows that describe each data sheet.
$keySheet = $wb.Worksheets.Item("Key")
$keySheet.Range("A3").Value = "cor_k07"
$keySheet.Range("B3").Value = "Correlation values vs. manual annotation for k07"
$keySheet.Range("A4").Value = "cor_k09"
$keySheet.Range("B4").Value = "Correlation values vs. manual annotation for k09"
